# Add "0x" prefix to each colon-separated hex byte in columns G (doip) and H (uds)
# for all data rows (rows 2..38), leaving "N/A" values untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($col in @("G", "H")) {
        $cell = $ws.Range("$col$r")
        $val = $cell.Value2
        if ($val -and $val -ne "N/A" -and $val -like "*:*") {
            $parts = $val -split ":"
            $newParts = @()
            foreach ($p in $parts) {
                $newParts += "0x$p"
            }
            $newVal = $newParts -join ":"
            $cell.Value = $newVal
        }
    }
}
